$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three Target-cluster = "ECs" rows (old rows 2, 5, 8), deleting
# from the bottom up so earlier row numbers stay valid.
$ws.Rows(8).EntireRow.Delete()
$ws.Rows(5).EntireRow.Delete()
$ws.Rows(2).EntireRow.Delete()

# Refresh the remaining six rows with the updated TPM-derived values.
$ws.Range("A2").Value = 'ECs'
$ws.Range("B2").Value = 'Gm13306'
$ws.Range("C2").Value = 'Ccr10'
$ws.Range("D2").Value = 'FAPs'
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.044816
$ws.Range("H2").Value = 0.134448
$ws.Range("I2").Value = 0.031000309200692
$ws.Range("J2").Value = 0.031000309200692
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.238415
$ws.Range("N2").Value = 3.715245
$ws.Range("O2").Value = 0.4359607654144799
$ws.Range("P2").Value = 0.4359607654144798
$ws.Range("Q2").Value = 0.05550080664
$ws.Range("R2").Value = 0.4995072597600001
$ws.Range("S2").Value = 0.01351491852721923
$ws.Range("T2").Value = 0.01351491852721923
$ws.Range("A3").Value = 'ECs'
$ws.Range("B3").Value = 'Gm13306'
$ws.Range("C3").Value = 'Ccr10'
$ws.Range("D3").Value = 'MuSCs'
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.044816
$ws.Range("H3").Value = 0.134448
$ws.Range("I3").Value = 0.031000309200692
$ws.Range("J3").Value = 0.031000309200692
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.602242
$ws.Range("N3").Value = 4.806725999999999
$ws.Range("O3").Value = 0.5640392345855201
$ws.Range("P3").Value = 0.5640392345855201
$ws.Range("Q3").Value = 0.07180607747199999
$ws.Range("R3").Value = 0.646254697248
$ws.Range("S3").Value = 0.01748539067347277
$ws.Range("T3").Value = 0.01748539067347277
$ws.Range("A4").Value = 'FAPs'
$ws.Range("B4").Value = 'Gm13306'
$ws.Range("C4").Value = 'Ccr10'
$ws.Range("D4").Value = 'FAPs'
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.8224453333333334
$ws.Range("H4").Value = 2.467336
$ws.Range("I4").Value = 0.5689052935112355
$ws.Range("J4").Value = 0.5689052935112355
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.238415
$ws.Range("N4").Value = 3.715245
$ws.Range("O4").Value = 0.4359607654144799
$ws.Range("P4").Value = 0.4359607654144798
$ws.Range("Q4").Value = 1.01852863748
$ws.Range("R4").Value = 9.166757737320001
$ws.Range("S4").Value = 0.2480203872075076
$ws.Range("T4").Value = 0.2480203872075075
$ws.Range("A5").Value = 'FAPs'
$ws.Range("B5").Value = 'Gm13306'
$ws.Range("C5").Value = 'Ccr10'
$ws.Range("D5").Value = 'MuSCs'
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.8224453333333334
$ws.Range("H5").Value = 2.467336
$ws.Range("I5").Value = 0.5689052935112355
$ws.Range("J5").Value = 0.5689052935112355
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.602242
$ws.Range("N5").Value = 4.806725999999999
$ws.Range("O5").Value = 0.5640392345855201
$ws.Range("P5").Value = 0.5640392345855201
$ws.Range("Q5").Value = 1.317756455770666
$ws.Range("R5").Value = 11.859808101936
$ws.Range("S5").Value = 0.320884906303728
$ws.Range("T5").Value = 0.320884906303728
$ws.Range("A6").Value = 'MuSCs'
$ws.Range("B6").Value = 'Gm13306'
$ws.Range("C6").Value = 'Ccr10'
$ws.Range("D6").Value = 'FAPs'
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5784016666666667
$ws.Range("H6").Value = 1.735205
$ws.Range("I6").Value = 0.4000943972880724
$ws.Range("J6").Value = 0.4000943972880724
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.238415
$ws.Range("N6").Value = 3.715245
$ws.Range("O6").Value = 0.4359607654144799
$ws.Range("P6").Value = 0.4359607654144798
$ws.Range("Q6").Value = 0.716301300025
$ws.Range("R6").Value = 6.446711700225001
$ws.Range("S6").Value = 0.1744254596797531
$ws.Range("T6").Value = 0.174425459679753
$ws.Range("A7").Value = 'MuSCs'
$ws.Range("B7").Value = 'Gm13306'
$ws.Range("C7").Value = 'Ccr10'
$ws.Range("D7").Value = 'MuSCs'
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5784016666666667
$ws.Range("H7").Value = 1.735205
$ws.Range("I7").Value = 0.4000943972880724
$ws.Range("J7").Value = 0.4000943972880724
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.602242
$ws.Range("N7").Value = 4.806725999999999
$ws.Range("O7").Value = 0.5640392345855201
$ws.Range("P7").Value = 0.5640392345855201
$ws.Range("Q7").Value = 0.9267394432033332
$ws.Range("R7").Value = 8.34065498883
$ws.Range("S7").Value = 0.2256689376083194
$ws.Range("T7").Value = 0.2256689376083194
